$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Hunk 0
$ws.Range("H4").Value = 544.6667
$ws.Range("I4").Value = 544.6667
$ws.Range("K4").Value = 544.6667
$ws.Range("M4").Value = -430.6667
# Hunk 1
$ws.Range("H18").Value = 111112420
$ws.Range("I18").Value = 111112420
$ws.Range("K18").Value = 111112420
$ws.Range("M18").Value = -111112136
# Hunk 2
$ws.Range("H40").Value = 55581730
$ws.Range("I40").Value = 27385
$ws.Range("K40").Value = 27385
$ws.Range("M40").Value = -27210
# Hunk 3
$ws.Range("H121").Value = 3427.56
$ws.Range("J121").Value = 3520.5833
$ws.Range("L121").Value = 10561.7499
$ws.Range("N121").Value = -14055.7499
# Hunk 4
$ws.Range("H132").Value = 103368.555
$ws.Range("I132").Value = 228722.42
$ws.Range("J132").Value = 13830.071
$ws.Range("K132").Value = 686167.26
$ws.Range("L132").Value = 41490.213
$ws.Range("M132").Value = -683637.26
$ws.Range("N132").Value = -46550.213
# Hunk 5
$ws.Range("H138").Value = 8131.5156
$ws.Range("I138").Value = 2234
$ws.Range("J138").Value = 9782.82
$ws.Range("K138").Value = 6702
$ws.Range("L138").Value = 29348.46
$ws.Range("M138").Value = -1562
$ws.Range("N138").Value = -39628.46
# Hunk 6
$ws.Range("H141").Value = 4685.44
$ws.Range("I141").Value = 4672.3335
$ws.Range("J141").Value = 5000
$ws.Range("K141").Value = 14017.0005
$ws.Range("L141").Value = 15000
$ws.Range("M141").Value = -8837.000499999998
$ws.Range("N141").Value = -25360

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Hunk 7
$ws.Range("H5").Value = 447.6
$ws.Range("I5").Value = 522
$ws.Range("J5").Value = 150
$ws.Range("K5").Value = 522
$ws.Range("L5").Value = 150
$ws.Range("M5").Value = -410
$ws.Range("N5").Value = -374
# Hunk 8
$ws.Range("H50").Value = 3313
$ws.Range("J50").Value = 3990.6667
$ws.Range("L50").Value = 3990.6667
$ws.Range("N50").Value = -5418.6667
# Hunk 9
$ws.Range("H74").Value = 3804.875
$ws.Range("I74").Value = 1425.3704
$ws.Range("J74").Value = 6864.2383
$ws.Range("K74").Value = 1425.3704
$ws.Range("L74").Value = 6864.2383
$ws.Range("M74").Value = -551.3704
$ws.Range("N74").Value = -8612.238300000001
# Hunk 10
$ws.Range("H77").Value = 3804.875
$ws.Range("I77").Value = 1425.3704
$ws.Range("J77").Value = 6864.2383
$ws.Range("K77").Value = 7126.852
$ws.Range("L77").Value = 34321.1915
$ws.Range("M77").Value = -2758.852
$ws.Range("N77").Value = -43057.1915

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Hunk 11
$ws.Range("H4").Value = 447.6
$ws.Range("I4").Value = 522
$ws.Range("J4").Value = 150
$ws.Range("K4").Value = 522
$ws.Range("L4").Value = 150
$ws.Range("M4").Value = -407
$ws.Range("N4").Value = -380
# Hunk 12
$ws.Range("H22").Value = 574.6667
$ws.Range("I22").Value = 574.6667
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 574.6667
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -401.6667
$ws.Range("N22").ClearContents()
# Hunk 13
$ws.Range("H140").Value = 163056.12
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 163056.12
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 163056.12
$ws.Range("M140").ClearContents()
$ws.Range("N140").Value = -173416.12

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Hunk 14
$ws.Range("H31").Value = 2301.08
$ws.Range("I31").Value = 2281.889
$ws.Range("J31").Value = 2323.6086
$ws.Range("K31").Value = 2281.889
$ws.Range("L31").Value = 2323.6086
$ws.Range("M31").Value = -1986.889
$ws.Range("N31").Value = -2913.6086
# Hunk 15
$ws.Range("H34").Value = 2301.08
$ws.Range("I34").Value = 2281.889
$ws.Range("J34").Value = 2323.6086
$ws.Range("K34").Value = 2281.889
$ws.Range("L34").Value = 2323.6086
$ws.Range("M34").Value = -2079.889
$ws.Range("N34").Value = -2727.6086
# Hunk 16
$ws.Range("H58").Value = 835041.75
$ws.Range("I58").Value = 1112389.1
$ws.Range("K58").Value = 1112389.1
$ws.Range("M58").Value = -1112186.1
# Hunk 17
$ws.Range("H59").Value = 95000
$ws.Range("J59").Value = 95000
$ws.Range("L59").Value = 95000
$ws.Range("N59").Value = -97290
# Hunk 18
$ws.Range("H60").Value = 233333
$ws.Range("J60").Value = 233333
$ws.Range("L60").Value = 233333
$ws.Range("N60").Value = -234355
# Hunk 19
$ws.Range("H99").Value = 10950.125
$ws.Range("I99").Value = 8567.333000000001
$ws.Range("J99").Value = 12379.8
$ws.Range("K99").Value = 8567.333000000001
$ws.Range("L99").Value = 12379.8
$ws.Range("M99").Value = -7069.333000000001
$ws.Range("N99").Value = -15375.8
# Hunk 20
$ws.Range("H123").Value = 42975
$ws.Range("I123").Value = 25000
$ws.Range("J123").Value = 48966.668
$ws.Range("K123").Value = 25000
$ws.Range("L123").Value = 48966.668
$ws.Range("M123").Value = -20100
$ws.Range("N123").Value = -58766.668
# Hunk 21
$ws.Range("H126").Value = 10950.125
$ws.Range("I126").Value = 8567.333000000001
$ws.Range("J126").Value = 12379.8
$ws.Range("K126").Value = 25701.999
$ws.Range("L126").Value = 37139.39999999999
$ws.Range("M126").Value = -23231.999
$ws.Range("N126").Value = -42079.39999999999
# Hunk 22
$ws.Range("H134").Value = 2520.9
$ws.Range("I134").Value = 2452
$ws.Range("J134").Value = 2911.3333
$ws.Range("K134").Value = 7356
$ws.Range("L134").Value = 8733.999899999999
$ws.Range("M134").Value = -4821
$ws.Range("N134").Value = -13803.9999
# Hunk 23
$ws.Range("H136").Value = 835041.75
$ws.Range("I136").Value = 1112389.1
$ws.Range("K136").Value = 3337167.3
$ws.Range("M136").Value = -3334617.3

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Hunk 24
$ws.Range("H2").Value = 2295844
$ws.Range("I2").Value = 4208796.5
$ws.Range("J2").Value = 300.8
$ws.Range("K2").Value = 25252779
$ws.Range("L2").Value = 1804.8
$ws.Range("M2").Value = -25252666
$ws.Range("N2").Value = -2030.8
# Hunk 25
$ws.Range("H38").Value = 146.73685
$ws.Range("I38").Value = 53
$ws.Range("J38").Value = 275.625
$ws.Range("K38").Value = 159
$ws.Range("L38").Value = 826.875
$ws.Range("M38").Value = 188
$ws.Range("N38").Value = -1520.875
# Hunk 26
$ws.Range("H107").Value = 3579.4
$ws.Range("I107").Value = 1974.25
$ws.Range("J107").Value = 10000
$ws.Range("K107").Value = 5922.75
$ws.Range("L107").Value = 30000
$ws.Range("M107").Value = -4002.75
$ws.Range("N107").Value = -33840
# Hunk 27
$ws.Range("H128").Value = 344965
$ws.Range("I128").Value = 344965
$ws.Range("K128").Value = 1034895
$ws.Range("M128").Value = -1029915

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Hunk 28
$ws.Range("H46").Value = 4687.5
$ws.Range("I46").Value = 1618.8889
$ws.Range("J46").Value = 6141.0527
$ws.Range("K46").Value = 1618.8889
$ws.Range("L46").Value = 6141.0527
$ws.Range("M46").Value = -1430.8889
$ws.Range("N46").Value = -6517.0527
# Hunk 29
$ws.Range("H132").Value = 4313.15
$ws.Range("I132").Value = 4313.15
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 12939.45
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -10409.45
$ws.Range("N132").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Hunk 30
$ws.Range("H136").Value = 9380.040000000001
$ws.Range("I136").Value = 6058.143
$ws.Range("J136").Value = 9920.813
$ws.Range("K136").Value = 18174.429
$ws.Range("L136").Value = 29762.439
$ws.Range("M136").Value = -15624.429
$ws.Range("N136").Value = -34862.439
